$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; temporarily unprotect so the refreshed
# holdings data (and the "as of" disclosure date) can be written.
$ws.Unprotect("D382")

# Update the "as of" date in the confidential disclosure note (cell A38),
# using an in-place text replace so the shared-string entry is edited
# rather than duplicated.
$null = $ws.Cells.Replace("2021-04-27", "2021-04-28", 2, 1, $false, $false, $false, $false)

# Refresh the Weight (D) and Percent Change (E) figures for each holding (rows 2-35)
$ws.Range("D2").Value2 = 0.03632513317591172
$ws.Range("E2").Value2 = -0.001160092807424684
$ws.Range("D3").Value2 = 0.02051693170976713
$ws.Range("E3").Value2 = -0.001666155298102212
$ws.Range("D4").Value2 = 0.0192549000942019
$ws.Range("E4").Value2 = 0.002016129032258229
$ws.Range("D5").Value2 = 0.03796937987060707
$ws.Range("E5").Value2 = 0.0003497726477790231
$ws.Range("D6").Value2 = 0.03528525823926655
$ws.Range("E6").Value2 = 0.002336448598130758
$ws.Range("D7").Value2 = 0.01985968048765307
$ws.Range("E7").Value2 = 0.001929012345678993
$ws.Range("D8").Value2 = 0.03657797633026841
$ws.Range("E8").Value2 = 0.005627705627705648
$ws.Range("D9").Value2 = 0.02028629584627672
$ws.Range("E9").Value2 = 0.003353879622915334
$ws.Range("D10").Value2 = 0.02518471055668041
$ws.Range("E10").Value2 = -0.0004056383733901381
$ws.Range("D11").Value2 = 0.02361784919104574
$ws.Range("E11").Value2 = -0.0002703433360368157
$ws.Range("D12").Value2 = 0.05700770320228929
$ws.Range("E12").Value2 = -0.00311004784688973
$ws.Range("D13").Value2 = 0.02508204091218406
$ws.Range("E13").Value2 = 0.0003665689149559004
$ws.Range("D14").Value2 = 0.02694470322870349
$ws.Range("E14").Value2 = 0.0001554484688326419
$ws.Range("D15").Value2 = 0.0322200840965117
$ws.Range("E15").Value2 = 0.001953125
$ws.Range("D16").Value2 = 0.01905083779083726
$ws.Range("E16").Value2 = 0.009362808842652726
$ws.Range("D17").Value2 = 0.03054805019454968
$ws.Range("E17").Value2 = -0.007022824178580356
$ws.Range("D18").Value2 = 0.042469477223905
$ws.Range("E18").Value2 = 0.0002285191956126198
$ws.Range("D19").Value2 = 0.1266473482373888
$ws.Range("E19").Value2 = 0.0006622516556291647
$ws.Range("D20").Value2 = 0.00908845995320765
$ws.Range("E20").Value2 = -0.003490175801447615
$ws.Range("D21").Value2 = 0.0154131654513087
$ws.Range("E21").Value2 = -0.001322291043218016
$ws.Range("D22").Value2 = 0.01695974828516971
$ws.Range("E22").Value2 = -0.02027545952588194
$ws.Range("D23").Value2 = 0.016012582505122
$ws.Range("E23").Value2 = -0.001738525730180851
$ws.Range("D24").Value2 = 0.02161052994257921
$ws.Range("E24").Value2 = -0.004167090151438124
$ws.Range("D25").Value2 = 0.01199523679865743
$ws.Range("E25").Value2 = 0.02368045649072759
$ws.Range("D26").Value2 = 0.04131916857014484
$ws.Range("E26").Value2 = -0.0001112594570538139
$ws.Range("D27").Value2 = 0.02395001869261779
$ws.Range("E27").Value2 = 0.0000981065437062334
$ws.Range("D28").Value2 = 0.04593501190077051
$ws.Range("E28").Value2 = 0.0028422548555187
$ws.Range("D29").Value2 = 0.05624498522635917
$ws.Range("E29").Value2 = -0.00605197579209682
$ws.Range("D30").Value2 = 0.01308338179204342
$ws.Range("E30").Value2 = 0.01555411535968876
$ws.Range("D31").Value2 = 0.02063251218975746
$ws.Range("E31").Value2 = 0.001534919416730673
$ws.Range("D32").Value2 = 0.01434054860315997
$ws.Range("E32").Value2 = 0.001780943900267173
$ws.Range("D33").Value2 = 0.04185197805359832
$ws.Range("E33").Value2 = 0.00103252452245739
$ws.Range("D34").Value2 = 0.01671431164745579
$ws.Range("E34").Value2 = 0.003144654088050203
$ws.Range("D35").Value2 = 0.9999999999999999
$ws.Range("E35").Value2 = 0.0002248393703276896

# Re-apply the original sheet protection.
$ws.Protect("D382")
